$doc = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the row for n = 10 in the "Mesures" table (2nd table),
#    since it always measures 0 seconds and is not relevant.
# ------------------------------------------------------------------
$table = $doc.Tables.Item(2)
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $row = $table.Rows.Item($i)
    $firstCellText = $row.Cells.Item(1).Range.Text
    $firstCellText = $firstCellText.TrimEnd([char]13, [char]7)
    if ($firstCellText -eq "10") {
        $row.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the paragraph that holds the
#    introductory image to the (now) very last, empty paragraph of
#    the document.
# ------------------------------------------------------------------
$goBack = $doc.Bookmarks.Item("_GoBack")
$goBack.Delete()

$lastPara = $doc.Paragraphs.Last
$lastRange = $lastPara.Range
# Temporarily insert a character so the target position is no longer
# the absolute last character of the document (avoids the COM host's
# end-of-story range quirk), add the bookmark, then remove the
# character again, leaving an empty paragraph with the bookmark.
$lastRange.InsertBefore("X")
$lastPara2 = $doc.Paragraphs.Last
$insertPos = $doc.Range($lastPara2.Range.Start, $lastPara2.Range.Start)
$doc.Bookmarks.Add("_GoBack", $insertPos)
$lastPara3 = $doc.Paragraphs.Last
$tempChar = $doc.Range($lastPara3.Range.Start, $lastPara3.Range.Start + 1)
$tempChar.Delete()
